$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Description" text in C5 (step list), incrementing the counts
# mentioned in steps 2 and 3.
$ws.Range("C5").Value = "1) All the blocks in the buffer pool (except log manager) are pinned once in order 11 to 17
2) Block 17, Block 15 and Block 13 and pinned 3 times
3) Block 17, Block 15 and Block 13 are unpinned 4 times
4) Block 18 and Block 19 will be newly pinned"

# Update the "Comments" text in E5, incrementing the pin/unpin counts.
$ws.Range("E5").Value = "Block 13, Block 15 and Block 17 were pinned 4 times (int total) and unpinned 4 times to ensure that the pin count remains zero at the end.
Block 15 will be replaced because Backward 2 distance of Block 17 is higher than Block 15 (as the second pin of Block 17 was before Block 15)"

# Move the active selection to E5, as left by the editor after making the change.
$ws.Range("E5").Select()
